$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, shifting existing rows 110:155 down to 111:156.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new data record.
$ws.Cells.Item(110, 1).Value = 4
$ws.Cells.Item(110, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(110, 3).Value = "Los Lagos"
$ws.Cells.Item(110, 4).Value = 44553
$ws.Cells.Item(110, 5).Value = 10
$ws.Cells.Item(110, 6).Value = 100112039
$ws.Cells.Item(110, 7).Value = "Ciboulette"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 120
$ws.Cells.Item(110, 11).Value = 2500
$ws.Cells.Item(110, 12).Value = 2500
$ws.Cells.Item(110, 13).Value = 2500
$ws.Cells.Item(110, 14).Value = "`$/docena de atados"
$ws.Cells.Item(110, 15).Value = "Región Metropolitana"
$ws.Cells.Item(110, 16).Value = 833
$ws.Cells.Item(110, 17).Value = 3
$ws.Cells.Item(110, 18).Value = "Hortaliza"

# Match the D-column (Fecha) number/date style used throughout the table.
$ws.Cells.Item(110, 4).NumberFormat = $ws.Cells.Item(111, 4).NumberFormat
